$wb = $excel.ActiveWorkbook

# --- Rename sheets ---
# rId1 sheet: "Q4_19_20" -> "Q1_20_21"
# rId2 sheet: "Q4_18_19" -> "Q4_19_20"
# Use a temp name on sheet1 first to avoid a name collision with sheet2 target name.
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)
$ws1.Name = "TEMP_RENAME_1"
$ws2.Name = "Q4_19_20"
$ws1.Name = "Q1_20_21"

# --- Sheet1 (rId1) data changes ---
$ws1.Range("B3").Value = "SRO"
$ws1.Range("C4").Value = 1
$ws1.Range("D4").Value = 469.7
$ws1.Range("E4").Value = 0.02622854943985702
$ws1.Range("C5").Value = 3
$ws1.Range("D5").Value = 10686.66708285521
$ws1.Range("E5").Value = 0.5967548987225048
$ws1.Range("C6").Value = 1
$ws1.Range("D6").Value = 6751.6
$ws1.Range("E6").Value = 0.3770165518376383
$ws1.Range("C7").Value = 0
$ws1.Range("D7").Value = 0
$ws1.Range("E7").Value = 0
$ws1.Range("C10").Value = 5
$ws1.Range("D10").Value = 17907.96708285521
$ws1.Range("B13").Value = "Finance"
$ws1.Range("C14").Value = 1
$ws1.Range("D14").Value = 469.7
$ws1.Range("E14").Value = 0.02622854943985702
$ws1.Range("C15").Value = 2
$ws1.Range("D15").Value = 7964.63
$ws1.Range("E15").Value = 0.4447534420378293
$ws1.Range("C16").Value = 2
$ws1.Range("D16").Value = 9473.63708285521
$ws1.Range("E16").Value = 0.5290180085223137
$ws1.Range("C20").Value = 5
$ws1.Range("D20").Value = 17907.96708285521
$ws1.Range("B23").Value = "Benefits"
$ws1.Range("C24").Value = 1
$ws1.Range("D24").Value = 1478.68
$ws1.Range("E24").Value = 0.08257106980142173
$ws1.Range("C25").Value = 2
$ws1.Range("D25").Value = 9207.98708285521
$ws1.Range("E25").Value = 0.5141838289210832
$ws1.Range("C26").Value = 2
$ws1.Range("D26").Value = 7221.3
$ws1.Range("E26").Value = 0.4032451012774953
$ws1.Range("C30").Value = 5
$ws1.Range("D30").Value = 17907.96708285521
$ws1.Range("B33").Value = "Schedule"
$ws1.Range("C34").Value = 4
$ws1.Range("D34").Value = 11156.36708285521
$ws1.Range("E34").Value = 0.6229834481623618
$ws1.Range("C36").Value = 0
$ws1.Range("D36").Value = 0
$ws1.Range("E36").Value = 0
$ws1.Range("C37").Value = 1
$ws1.Range("D37").Value = 6751.6
$ws1.Range("E37").Value = 0.3770165518376383
$ws1.Range("C38").Value = 0
$ws1.Range("D38").Value = 0
$ws1.Range("E38").Value = 0
$ws1.Range("C40").Value = 5
$ws1.Range("D40").Value = 17907.96708285521
$ws1.Range("B43").Value = "Resource"
$ws1.Range("C44").Value = 2
$ws1.Range("D44").Value = 7221.3
$ws1.Range("E44").Value = 0.4032451012774953
$ws1.Range("C45").Value = 3
$ws1.Range("D45").Value = 10686.66708285521
$ws1.Range("E45").Value = 0.5967548987225048
$ws1.Range("C49").Value = 0
$ws1.Range("D49").Value = 0
$ws1.Range("E49").Value = 0
$ws1.Range("C50").Value = 5
$ws1.Range("D50").Value = 17907.96708285521

# --- Sheet2 (rId2) data changes ---
$ws2.Range("B3").Value = "SRO"
$ws2.Range("C4").Value = 2
$ws2.Range("D4").Value = 2118.6
$ws2.Range("E4").Value = 0.1202439341118181
$ws2.Range("C5").Value = 1
$ws2.Range("D5").Value = 6490.37
$ws2.Range("E5").Value = 0.3683694999723029
$ws2.Range("C6").Value = 2
$ws2.Range("D6").Value = 9010.214
$ws2.Range("E6").Value = 0.5113865659158791
$ws2.Range("D10").Value = 17619.184
$ws2.Range("B13").Value = "Finance"
$ws2.Range("C14").Value = 2
$ws2.Range("D14").Value = 2118.6
$ws2.Range("E14").Value = 0.1202439341118181
$ws2.Range("C16").Value = 3
$ws2.Range("D16").Value = 15500.584
$ws2.Range("E16").Value = 0.8797560658881819
$ws2.Range("D20").Value = 17619.184
$ws2.Range("B23").Value = "Benefits"
$ws2.Range("C24").Value = 2
$ws2.Range("D24").Value = 4064.784
$ws2.Range("E24").Value = 0.2307021709972494
$ws2.Range("C26").Value = 3
$ws2.Range("D26").Value = 13554.4
$ws2.Range("E26").Value = 0.7692978290027507
$ws2.Range("C29").Value = 0
$ws2.Range("D29").Value = 0
$ws2.Range("E29").Value = 0
$ws2.Range("D30").Value = 17619.184
$ws2.Range("B33").Value = "Schedule"
$ws2.Range("D34").Value = 2118.6
$ws2.Range("E34").Value = 0.1202439341118181
$ws2.Range("C36").Value = 2
$ws2.Range("D36").Value = 9108.853999999999
$ws2.Range("E36").Value = 0.5169850090673893
$ws2.Range("C38").Value = 1
$ws2.Range("D38").Value = 6391.73
$ws2.Range("E38").Value = 0.3627710568207926
$ws2.Range("D40").Value = 17619.184
$ws2.Range("B43").Value = "Resource"
$ws2.Range("C43").Value = "Count"
$ws2.Range("D43").Value = "Costs"
$ws2.Range("E43").Value = "Proportion costs"
$ws2.Range("B44").Value = "Green"
$ws2.Range("C44").Value = 5
$ws2.Range("D44").Value = 17619.184
$ws2.Range("E44").Value = 1
$ws2.Range("B45").Value = "Amber/Green"
$ws2.Range("C45").Value = 0
$ws2.Range("D45").Value = 0
$ws2.Range("E45").Value = 0
$ws2.Range("B46").Value = "Amber"
$ws2.Range("C46").Value = 0
$ws2.Range("D46").Value = 0
$ws2.Range("E46").Value = 0
$ws2.Range("B47").Value = "Amber/Red"
$ws2.Range("C47").Value = 0
$ws2.Range("D47").Value = 0
$ws2.Range("E47").Value = 0
$ws2.Range("B48").Value = "Red"
$ws2.Range("C48").Value = 0
$ws2.Range("D48").Value = 0
$ws2.Range("E48").Value = 0
$ws2.Range("B49").Value = "None"
$ws2.Range("C49").Value = 0
$ws2.Range("D49").Value = 0
$ws2.Range("E49").Value = 0
$ws2.Range("B50").Value = "Total"
$ws2.Range("C50").Value = 5
$ws2.Range("D50").Value = 17619.184
$ws2.Range("E50").Value = 1

